# Applies the changes described by the diff to pesquisas-infos.xlsx (Sheet1):
#  1. Fix typo "9-2014" -> "09-2014" in D48, D49, D50
#  2. Append 6 new data rows (111-116) for id_pesquisa 04843 / p901a question
#  3. Dimension will automatically extend to A1:L116 once the new cells are written

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, [string]$Address, [string]$Value)
    # Force the cell to be treated as plain text so that values which look
    # numeric (e.g. "04843", "2010", "23") keep leading zeros / stay strings
    # instead of being silently converted to numbers by Excel.
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value2 = $Value
    $cell.Style = "Normal"
}

# --- 1. Fix the "9-2014" -> "09-2014" typo -------------------------------
$ws.Range("D48").Value2 = "09-2014"
$ws.Range("D49").Value2 = "09-2014"
$ws.Range("D50").Value2 = "09-2014"

# --- 2. Append new rows 111-116 -------------------------------------------
$newRows = @(
    @{ Row = 111; J = "Não sabe";                     K = "23" },
    @{ Row = 112; J = "Concorda totalmente";           K = "602" },
    @{ Row = 113; J = "Concorda em parte";              K = "346" },
    @{ Row = 114; J = "Discorda em parte";               K = "220" },
    @{ Row = 115; J = "Nem concorda, nem discorda";       K = "16" },
    @{ Row = 116; J = "Discorda totalmente";               K = "803" }
)

$A = "04843"
$B = "Data Folha"
$C = "PERFIL IDEOLÓGICO"
$D = "06-2023"
$E = "2010"
$F = "Eleitores de 16 anos ou mais"
$G = "Brasil"
$H = "p901a"
$I = "P.901 Agora vou ler algumas frases e gostaria que você me dissesse se concorda ou discorda de cada uma delas. - Possuir arma legalizada deveria ser um direito do cidadão para se defender"
$L = "2010"

foreach ($r in $newRows) {
    $row = $r.Row

    Set-TextCell $ws "A$row" $A
    $ws.Range("B$row").Value2 = $B
    $ws.Range("C$row").Value2 = $C
    $ws.Range("D$row").Value2 = $D
    Set-TextCell $ws "E$row" $E
    $ws.Range("F$row").Value2 = $F
    $ws.Range("G$row").Value2 = $G
    $ws.Range("H$row").Value2 = $H
    $ws.Range("I$row").Value2 = $I
    $ws.Range("J$row").Value2 = $r.J
    Set-TextCell $ws "K$row" $r.K
    Set-TextCell $ws "L$row" $L
}
